$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(24, 24, "Constantine", "Ibn Ziad", "66677dd8d12a2002d45edd0a"),
    @(25, 24, "Constantine", "Ain Abid", "66677ea9d12a2002d45edd33"),
    @(26, 24, "Constantine", "Ibn Ziad", "66677f04d12a2002d45edd41"),
    @(27, 24, "Constantine", "El Khroub", "66677f75d12a2002d45edd54")
)

$startRow = 25
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
